$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws1.Range("A1").Value = "test"
